# Fix the "Инфоцентр" column (J) and row (10) which were left as zeros:
# - Column J (rows 2-9) should mirror column I's values for those rows.
# - Row 10 (columns B-I) should mirror row 9's values for those columns.
# Also apply the same number-style (thin-border/centered, style used by
# column I / row 9) to these previously-unstyled cells, and drop the
# stray "applyFill" flag from the bold corner-label style (A10 / J1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Column J: copy value + style from column I, rows 2-9 ---
for ($r = 2; $r -le 9; $r++) {
    $srcI = $ws.Cells.Item($r, 9)   # column I
    $dstJ = $ws.Cells.Item($r, 10)  # column J

    $srcI.Copy() | Out-Null
    $dstJ.PasteSpecial($xlPasteFormats) | Out-Null

    $dstJ.Value = $srcI.Value2
}

# --- Row 10: copy value + style from row 9, columns B-I ---
for ($c = 2; $c -le 9; $c++) {
    $src9 = $ws.Cells.Item(9, $c)   # row 9
    $dst10 = $ws.Cells.Item(10, $c) # row 10

    $src9.Copy() | Out-Null
    $dst10.PasteSpecial($xlPasteFormats) | Out-Null

    $dst10.Value = $src9.Value2
}

$excel.CutCopyMode = 0

# --- Remove the "applyFill" flag from the bold bordered style used by
#     the corner cells (A10 / J1) by clearing any explicit fill so the
#     cell no longer carries an applied-fill flag. ---
$ws.Range("A10").Interior.Pattern = -4142   # xlPatternNone
$ws.Range("J1").Interior.Pattern = -4142    # xlPatternNone
